$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text representation (avoid Excel
# auto-converting numeric-looking strings like "1.00" or "24.60" into
# true numbers and losing formatting / trailing zeros).
$cells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'E10', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'E15', 'D16', 'E16', 'D17', 'E17', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'E28', 'D29', 'E29', 'E30', 'E31', 'D32', 'E32', 'E33', 'D34', 'E34', 'E35', 'B36', 'C36', 'D36', 'E36', 'B37', 'C37', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'E41', 'E42', 'B43', 'C43', 'D43', 'E43', 'B44', 'C44', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.308.61'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '3.032.77'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '577.31'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = '168.07'
$ws.Range('E6').Value = '  +3.75%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.029.55'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '0.481'
$ws.Range('E12').Value = '  +5.59%  '
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').Value = '36.40'
$ws.Range('E14').Value = '  +5.37%  '
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '66.248.37'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '3.535.04'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D19').Value = '16.57'
$ws.Range('E19').Value = '  +20.02%  '
$ws.Range('D20').Value = '3.031.48'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = '474.91'
$ws.Range('E21').Value = '  +4.28%  '
$ws.Range('D22').Value = '0.709'
$ws.Range('E22').Value = '  +3.28%  '
$ws.Range('D23').Value = '7.46'
$ws.Range('E23').Value = '  +1.69%  '
$ws.Range('D24').Value = '83.24'
$ws.Range('E24').Value = '  +1.25%  '
$ws.Range('D25').Value = '12.86'
$ws.Range('E25').Value = '  +4.84%  '
$ws.Range('D26').Value = '2.27'
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('D27').Value = '10.07'
$ws.Range('E27').Value = '  -3.36%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = '8.22'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  +6.43%  '
$ws.Range('E33').Value = '  -5.88%  '
$ws.Range('D34').Value = '28.07'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '5.88'
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').Value = '0.991'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Value = '48.13'
$ws.Range('E38').Value = '  +9.75%  '
$ws.Range('D39').Value = '2.07'
$ws.Range('E39').Value = '  -3.95%  '
$ws.Range('D40').Value = '49.64'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '8.64'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '2.84'
$ws.Range('E44').Value = '  -4.87%  '
$ws.Range('D45').Value = '0.0361'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').Value = '384.09'
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('D47').Value = '2.724.78'
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('D48').Value = '134.74'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '24.60'
$ws.Range('E50').Value = '  +3.23%  '
$ws.Range('D51').Value = '2.24'
$ws.Range('E51').Value = '  +4.66%  '
